$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.04558966666666667
$ws.Range("H2").Value = 0.136769
$ws.Range("I2").Value = 0.02375599288687187
$ws.Range("J2").Value = 0.02375599288687187
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.040495
$ws.Range("N2").Value = 0.121485
$ws.Range("O2").Value = 0.002191743187342868
$ws.Range("P2").Value = 0.002191743187342869
$ws.Range("Q2").Value = 0.001846153551666666
$ws.Range("R2").Value = 0.016615381965
$ws.Range("S2").Value = 0.00005206703556836706
$ws.Range("T2").Value = 0.00005206703556836707

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.04558966666666667
$ws.Range("H3").Value = 0.136769
$ws.Range("I3").Value = 0.02375599288687187
$ws.Range("J3").Value = 0.02375599288687187
$ws.Range("O3").Value = 0.9191911494312409
$ws.Range("P3").Value = 0.9191911494312409
$ws.Range("Q3").Value = 0.7742549469221112
$ws.Range("R3").Value = 6.968294522299
$ws.Range("S3").Value = 0.02183629840756414
$ws.Range("T3").Value = 0.02183629840756414

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.04558966666666667
$ws.Range("H4").Value = 0.136769
$ws.Range("I4").Value = 0.02375599288687187
$ws.Range("J4").Value = 0.02375599288687187
$ws.Range("M4").Value = 1.452542333333333
$ws.Range("N4").Value = 4.357627
$ws.Range("O4").Value = 0.07861710738141615
$ws.Range("P4").Value = 0.07861710738141615
$ws.Range("Q4").Value = 0.06622092079588889
$ws.Range("R4").Value = 0.595988287163
$ws.Range("S4").Value = 0.001867627443739364
$ws.Range("T4").Value = 0.001867627443739364

# Row 5
$ws.Range("I5").Value = 0.1978186777627204
$ws.Range("J5").Value = 0.1978186777627204
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.040495
$ws.Range("N5").Value = 0.121485
$ws.Range("O5").Value = 0.002191743187342868
$ws.Range("P5").Value = 0.002191743187342869
$ws.Range("Q5").Value = 0.01537311685
$ws.Range("R5").Value = 0.13835805165
$ws.Range("S5").Value = 0.0004335677393156165
$ws.Range("T5").Value = 0.0004335677393156166

# Row 6
$ws.Range("I6").Value = 0.1978186777627204
$ws.Range("J6").Value = 0.1978186777627204
$ws.Range("O6").Value = 0.9191911494312409
$ws.Range("P6").Value = 0.9191911494312409
$ws.Range("S6").Value = 0.1818331777916832
$ws.Range("T6").Value = 0.1818331777916832

# Row 7
$ws.Range("I7").Value = 0.1978186777627204
$ws.Range("J7").Value = 0.1978186777627204
$ws.Range("M7").Value = 1.452542333333333
$ws.Range("N7").Value = 4.357627
$ws.Range("O7").Value = 0.07861710738141615
$ws.Range("P7").Value = 0.07861710738141615
$ws.Range("Q7").Value = 0.5514286460033333
$ws.Range("R7").Value = 4.962857814029999
$ws.Range("S7").Value = 0.01555193223172155
$ws.Range("T7").Value = 0.01555193223172155

# Row 8
$ws.Range("I8").Value = 0.7784253293504076
$ws.Range("J8").Value = 0.7784253293504078
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.040495
$ws.Range("N8").Value = 0.121485
$ws.Range("O8").Value = 0.002191743187342868
$ws.Range("P8").Value = 0.002191743187342869
$ws.Range("Q8").Value = 0.060493901195
$ws.Range("R8").Value = 0.544445110755
$ws.Range("S8").Value = 0.001706108412458884
$ws.Range("T8").Value = 0.001706108412458885

# Row 9
$ws.Range("I9").Value = 0.7784253293504076
$ws.Range("J9").Value = 0.7784253293504078
$ws.Range("O9").Value = 0.9191911494312409
$ws.Range("P9").Value = 0.9191911494312409
$ws.Range("S9").Value = 0.7155216732319934
$ws.Range("T9").Value = 0.7155216732319936

# Row 10
$ws.Range("I10").Value = 0.7784253293504076
$ws.Range("J10").Value = 0.7784253293504078
$ws.Range("M10").Value = 1.452542333333333
$ws.Range("N10").Value = 4.357627
$ws.Range("O10").Value = 0.07861710738141615
$ws.Range("P10").Value = 0.07861710738141615
$ws.Range("S10").Value = 0.06119754770595524
$ws.Range("T10").Value = 0.06119754770595524
